$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks so they can be re-added cleanly in the new row order
$ws.Hyperlinks.Delete()

# --- Column width adjustments (col B: 52 -> 55, col D: 30 -> 32) ---
$ws.Columns.Item(2).ColumnWidth = 54.16
$ws.Columns.Item(4).ColumnWidth = 31.16

# --- Row data (title/price/url/score/skills shift down as three new AI listings
#     are inserted near the top and six more are appended near the bottom) ---

# Row 2
$ws.Range("A2").Value = "2025-11-24 18:27:35"
$ws.Range("B2").Value = "【Python/AI/GAS 開発者・PM向け】「業務委託・再委託」の経験に関する30分インタビュー"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5434693"
$ws.Range("G2").Value = 530
$ws.Range("H2").Value = "🔥AI,Python ◆開発"

# Row 3
$ws.Range("A3").Value = "2025-11-24 18:27:35"
$ws.Range("B3").Value = "【AI×Web】建設業向け材料発注自動化SaaSのMVP開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5440400"
$ws.Range("G3").Value = 435
$ws.Range("H3").Value = "🔥AI,Ai ◆開発,自動化"

# Row 4
$ws.Range("A4").Value = "2025-11-24 18:27:35"
$ws.Range("B4").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("G4").Value = 385
$ws.Range("H4").Value = "🔥AI,Ai ◆効率化"

# Row 5
$ws.Range("A5").Value = "2025-11-24 18:27:35"
$ws.Range("B5").Value = "製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5419380"
$ws.Range("G5").Value = 298
$ws.Range("H5").Value = "🔥AI,Ai"

# Row 6
$ws.Range("A6").Value = "2025-11-24 18:27:35"
$ws.Range("B6").Value = "【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5440052"
$ws.Range("G6").Value = 158
$ws.Range("H6").Value = "◆自動化,スクレイピング ◇管理"

# Row 7
$ws.Range("A7").Value = "2025-11-24 18:27:35"
$ws.Range("B7").Value = "マッチングサイト開発エンジニア募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5440077"
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = "◆開発 ◇サイト"

# Row 8
$ws.Range("A8").Value = "2025-11-24 18:27:35"
$ws.Range("B8").Value = "【急募】掲示板サイト(爆サイ)自動書き込みソフト開発者募集"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5439484"
$ws.Range("G8").Value = 93
$ws.Range("H8").Value = "◆開発 ◇サイト"

# Row 9
$ws.Range("A9").Value = "2025-11-24 18:27:35"
$ws.Range("B9").Value = "【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5431107"
$ws.Range("G9").Value = 68
$ws.Range("H9").Value = "◆開発"

# Row 10
$ws.Range("A10").Value = "2025-11-24 18:27:35"
$ws.Range("B10").Value = "【出張対応】Accessシステム改修・CSV読込・MySQLクラウド化・PDFデータ調整【急募】"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5440318"
$ws.Range("G10").Value = 48
$ws.Range("H10").Value = "◇MySQL"

# Row 11
$ws.Range("A11").Value = "2025-11-24 18:27:35"
$ws.Range("B11").Value = "【急募】シティヘブンの出勤情報を自動取得・管理したい!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5440436"
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = "◇管理"

# Row 12
$ws.Range("A12").Value = "2025-11-24 18:27:35"
$ws.Range("B12").Value = "【急募】各種ホームページ・Wordpressの修正・保守対応をお任せ!"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5440417"
$ws.Range("G12").Value = 25
$ws.Range("H12").Value = "○WordPress"

# Row 13
$ws.Range("A13").Value = "2025-11-24 18:27:35"
$ws.Range("B13").Value = "初回 n8n+Gemini+Typefully+GoogleスプレッドのX/Threads自動投稿システム"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5440440"
$ws.Range("G13").Value = 33
$ws.Range("H13").ClearContents()

# Row 14
$ws.Range("A14").Value = "2025-11-24 18:27:35"
$ws.Range("B14").Value = "急募 限定公開 PR 限定公開の仕事"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5440230"
$ws.Range("G14").Value = 25
$ws.Range("H14").ClearContents()

# Row 15
$ws.Range("A15").Value = "2025-11-24 18:27:35"
$ws.Range("B15").Value = "【急募】貸別荘収支表自動集計システム構築の依頼"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5440042"
$ws.Range("G15").Value = 25
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = "2025-11-24 18:27:35"
$ws.Range("B16").Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5440002"
$ws.Range("G16").Value = 13
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = "2025-11-24 18:27:35"
$ws.Range("B17").Value = "【Amazon出品・Excel】ブラウズノード設定/フラットファイル検証に詳しい方を募集"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5440204"
$ws.Range("G17").Value = 13
$ws.Range("H17").ClearContents()

# Row 18
$ws.Range("A18").Value = "2025-11-24 18:27:35"
$ws.Range("B18").Value = "Chatworkとn8n連携"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5440325"
$ws.Range("G18").Value = 10
$ws.Range("H18").ClearContents()

# --- Re-add hyperlinks for every URL cell (F3:F18 are new targets; F2 unchanged) ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5434693")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5440400")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5423720")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5419380")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5440052")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5440077")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5439484")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5431107")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5440318")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5440436")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5440417")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5440440")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5440230")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5440042")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5440002")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5440204")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5440325")
